# Apply updated "dSF" (column F) values to Sheet1, per the source diff.
# Only column F values change; all other columns (including E / dS0) stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    2  = -6
    3  = -7
    4  = -3
    5  = -1
    6  = 4
    7  = 11
    12 = -7
    14 = -1
    15 = -1
    16 = -4
    17 = -5
    18 = 12
    20 = -4
    21 = 4
    22 = -1
    23 = -2
    24 = -2
    25 = 3
    26 = -1
    27 = -4
    28 = 2
    29 = 2
    30 = 1
    31 = 5
    33 = 5
    34 = 7
    35 = 5
    36 = 3
    37 = 3
    38 = -2
    39 = 1
    40 = -1
    42 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
